$wb = $excel.ActiveWorkbook

# --- Rename sheets (merged "Jobs" + "Papers" dynamic creation also renamed the
#     generic *Items sheets to their short display names) ---
$art = $wb.Worksheets.Item("ArtItems")
$art.Name = "Art"

$programming = $wb.Worksheets.Item("ProgrammingItems")
$programming.Name = "Programming"

# --- Jobs sheet: move the selection to A3 and make sure it is no longer the
#     tab that is active/selected (that moves to Programming below) ---
$jobs = $wb.Worksheets.Item("Jobs")
$jobs.Range("A3").Select()

# --- Programming sheet: add a new (formatted-then-cleared) cell in G21,
#     growing the used range, then activate the sheet and leave the
#     selection on G21 ---
$cell = $programming.Range("G21")
$cell.Font.Underline = $true
$cell.Font.Underline = $false

$programming.Activate()
$programming.Range("G21").Select()
